# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" sheet and the corresponding rows in the "全部类型" sheet.
# These mirror each other except the "全部类型" sheet has extra rows
# merged in from other categories, shifting the last entry's row number.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1415
$ws1.Range("F4").Value  = 19881
$ws1.Range("F9").Value  = 7511
$ws1.Range("F26").Value = 320
$ws1.Range("F27").Value = 1097
$ws1.Range("F31").Value = 5221
$ws1.Range("F41").Value = 23

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1415
$ws4.Range("F4").Value  = 19881
$ws4.Range("F9").Value  = 7511
$ws4.Range("F26").Value = 320
$ws4.Range("F27").Value = 1097
$ws4.Range("F31").Value = 5221
$ws4.Range("F43").Value = 23
